# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# for each coin row on Sheet1, per the scraped crypto-ranking snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is one changed cell (A1 ref) and its new text value, taken from
# the updated coinranking.com snapshot. Values must land as literal text
# (matching the workbook's existing inline-string cells), not as numbers.
$updates = @(
    @{ Cell = 'D2'; Value = '26.217.16' },
    @{ Cell = 'E2'; Value = '  -0.76%  ' },
    @{ Cell = 'D3'; Value = '1.672.00' },
    @{ Cell = 'E3'; Value = '  -1.44%  ' },
    @{ Cell = 'E4'; Value = '  -0.67%  ' },
    @{ Cell = 'D5'; Value = '211.67' },
    @{ Cell = 'E5'; Value = '  -3.17%  ' },
    @{ Cell = 'D6'; Value = '0.5281' },
    @{ Cell = 'E6'; Value = '  -3.77%  ' },
    @{ Cell = 'E7'; Value = '  -0.67%  ' },
    @{ Cell = 'D8'; Value = '0.2644' },
    @{ Cell = 'E8'; Value = '  -3.29%  ' },
    @{ Cell = 'D9'; Value = '0.06281' },
    @{ Cell = 'E9'; Value = '  -2.54%  ' },
    @{ Cell = 'D10'; Value = '21.31' },
    @{ Cell = 'E10'; Value = '  -2.94%  ' },
    @{ Cell = 'D11'; Value = '0.07564' },
    @{ Cell = 'E11'; Value = '  -1.54%  ' },
    @{ Cell = 'D12'; Value = '1.715.70' },
    @{ Cell = 'E12'; Value = '  +0.56%  ' },
    @{ Cell = 'D13'; Value = '4.459' },
    @{ Cell = 'E13'; Value = '  -2.14%  ' },
    @{ Cell = 'D14'; Value = '0.5596' },
    @{ Cell = 'E14'; Value = '  -4.32%  ' },
    @{ Cell = 'D15'; Value = '66.91' },
    @{ Cell = 'E15'; Value = '  +1.88%  ' },
    @{ Cell = 'D16'; Value = '0.000008016' },
    @{ Cell = 'E16'; Value = '  -4.61%  ' },
    @{ Cell = 'D17'; Value = '26.247.96' },
    @{ Cell = 'E17'; Value = '  -0.95%  ' },
    @{ Cell = 'E18'; Value = '  -0.65%  ' },
    @{ Cell = 'D19'; Value = '4.803' },
    @{ Cell = 'E19'; Value = '  -2.90%  ' },
    @{ Cell = 'D20'; Value = '187.53' },
    @{ Cell = 'E20'; Value = '  -2.09%  ' },
    @{ Cell = 'D21'; Value = '10.41' },
    @{ Cell = 'E21'; Value = '  -5.27%  ' },
    @{ Cell = 'D22'; Value = '6.209' },
    @{ Cell = 'E22'; Value = '  -0.79%  ' },
    @{ Cell = 'E23'; Value = '  -0.63%  ' },
    @{ Cell = 'D24'; Value = '150.02' },
    @{ Cell = 'E24'; Value = '  +0.76%  ' },
    @{ Cell = 'D25'; Value = '0.1258' },
    @{ Cell = 'E25'; Value = '  -4.35%  ' },
    @{ Cell = 'D26'; Value = '7.577' },
    @{ Cell = 'E26'; Value = '  -4.28%  ' },
    @{ Cell = 'D27'; Value = '15.97' },
    @{ Cell = 'E27'; Value = '  +0.80%  ' },
    @{ Cell = 'D28'; Value = '0.06217' },
    @{ Cell = 'E28'; Value = '  -0.32%  ' },
    @{ Cell = 'D29'; Value = '1.364' },
    @{ Cell = 'E29'; Value = '  -1.46%  ' },
    @{ Cell = 'E30'; Value = '  -3.58%  ' },
    @{ Cell = 'D31'; Value = '3.503' },
    @{ Cell = 'E31'; Value = '  -2.99%  ' },
    @{ Cell = 'D32'; Value = '3.430' },
    @{ Cell = 'E32'; Value = '  -4.64%  ' },
    @{ Cell = 'E33'; Value = '  -3.40%  ' },
    @{ Cell = 'E34'; Value = '  -3.64%  ' },
    @{ Cell = 'D35'; Value = '0.6048' },
    @{ Cell = 'E35'; Value = '  -1.86%  ' },
    @{ Cell = 'D36'; Value = '2.412' },
    @{ Cell = 'E36'; Value = '  +0.07%  ' },
    @{ Cell = 'D37'; Value = '2.744' },
    @{ Cell = 'E37'; Value = '  -0.66%  ' },
    @{ Cell = 'E38'; Value = '  +0.06%  ' },
    @{ Cell = 'D39'; Value = '0.01622' },
    @{ Cell = 'E39'; Value = '  -1.78%  ' },
    @{ Cell = 'D40'; Value = '1.102.73' },
    @{ Cell = 'E40'; Value = '  -1.47%  ' },
    @{ Cell = 'D41'; Value = '0.8760' },
    @{ Cell = 'E41'; Value = '  -0.90%  ' },
    @{ Cell = 'E42'; Value = '  -0.96%  ' },
    @{ Cell = 'D43'; Value = '99.95' },
    @{ Cell = 'E43'; Value = '  -1.27%  ' },
    @{ Cell = 'D44'; Value = '1.823.98' },
    @{ Cell = 'E44'; Value = '  -1.33%  ' },
    @{ Cell = 'D45'; Value = '0.00000000111' },
    @{ Cell = 'E45'; Value = '  +2.22%  ' },
    @{ Cell = 'D46'; Value = '55.99' },
    @{ Cell = 'E46'; Value = '  -2.86%  ' },
    @{ Cell = 'E47'; Value = '  -0.06%  ' },
    @{ Cell = 'D48'; Value = '8.020' },
    @{ Cell = 'E48'; Value = '  -2.55%  ' },
    @{ Cell = 'D49'; Value = '0.05229' },
    @{ Cell = 'E49'; Value = '  -1.12%  ' },
    @{ Cell = 'D50'; Value = '0.4255' },
    @{ Cell = 'E50'; Value = '  -1.17%  ' },
    @{ Cell = 'D51'; Value = '5.988' },
    @{ Cell = 'E51'; Value = '  -2.29%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $value = $u.Value

    # Plain decimal-looking price strings (e.g. "4.459", "0.000008016")
    # get auto-coerced to a real Number by COM's normal Value assignment,
    # which would corrupt the text formatting (and precision) that the
    # original sheet relies on. Force those to stay text: write with a
    # leading apostrophe (Excel's "treat as text" prefix), then put the
    # cell's style back to the workbook default so no quote-prefix /
    # text-number-format style is left behind on the cell.
    $looksNumeric = $value -match '^[+-]?\d+(\.\d+)?([eE][+-]?\d+)?$'

    if ($looksNumeric) {
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
